$wb = $excel.ActiveWorkbook

$wsFlows    = $wb.Worksheets.Item("_set_FLOWS")
$wsFlowsAgg = $wb.Worksheets.Item("_set_FLOWS_AGG")

# --- _set_FLOWS sheet : column widths -----------------------------------
# Collapse the three bestFit columns (A:C) into a single, explicit,
# user-set width so the three separate <col> entries become one uniform
# width definition (bestFit cleared, customWidth retained).
$wsFlows.Columns("A:C").ColumnWidth = 19.3

# --- _set_FLOWS_AGG sheet : reorder the aggregation rows -----------------
# Re-map which f_Name each aggregation row refers to, so that row order
# no longer matches the row order used in _set_FLOWS (this is the actual
# content fix described in the commit message).
$wsFlowsAgg.Range("A2").Value = "transport"
$wsFlowsAgg.Range("B2").Value = "Yearly dispatched"

$wsFlowsAgg.Range("A3").Value = "oil products"
$wsFlowsAgg.Range("B3").Value = "Yearly dispatched"

$wsFlowsAgg.Range("A4").Value = "electricity"
$wsFlowsAgg.Range("B4").Value = "Hourly dispatched"

# --- restore the selected cell on each sheet ------------------------------
$wsFlows.Select()
$wsFlows.Range("C14").Select() | Out-Null

$wsFlowsAgg.Select()
$wsFlowsAgg.Range("B9").Select() | Out-Null
